$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new request row (row 3) below the existing header (row 1) and
# data row (row 2), mirroring the same column layout.
$ws.Range("A3").Value = "20250527_095356"
$ws.Range("B3").Value = "2025-05-27 09:53:56"
$ws.Range("C3").Value = "Marie Martin"
$ws.Range("D3").Value = "{'chantier': 'Aluminium - Table Aluminium 02', 'urgence': 'Normal', 'date_souhaitee': '2025-05-27', 'produits': {'3813399991': {'produit': 'CYLINDRE RX 60 30X30 NICK MAT VARIE ', 'quantite': 1, 'emplacement': 'A11'}}}"
$ws.Range("F3").Value = "En attente"

# Motif / Date_Traitement / Traite_Par / Commentaires stay blank for a new
# request, same as the blank cells already present on row 2. Touch each
# cell's formatting so it is still materialized as an (empty) cell rather
# than being omitted entirely from the sheet.
$ws.Range("E3").Value = ""
$ws.Range("E3").Interior.Pattern = -4142
$ws.Range("G3").Value = ""
$ws.Range("G3").Interior.Pattern = -4142
$ws.Range("H3").Value = ""
$ws.Range("H3").Interior.Pattern = -4142
$ws.Range("I3").Value = ""
$ws.Range("I3").Interior.Pattern = -4142
